# "removed pally's hp buff"
#
# The Paladin class bullet list currently reads:
#   Paladin
#     - +2 Hit Points
#     - 1 Color Slot
#     - Create White
#     - ...
#
# This script removes the "+2 Hit Points" bullet point entirely (including
# its paragraph mark), leaving "Paladin" followed directly by "1 Color Slot".

$d = $word.ActiveDocument

# "Paladin" (the class heading run) appears exactly once in the whole
# document, so find its paragraph index by scanning the Paragraphs
# collection directly (keeps everything in document-level indices).
$headingIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.TrimEnd() -eq "Paladin") {
        $headingIndex = $i
    }
}

if ($headingIndex -eq -1) {
    throw "Could not find the 'Paladin' class heading."
}

# The very next paragraph is the "+2 Hit Points" bullet we need to drop.
$hpParagraph = $d.Paragraphs.Item($headingIndex + 1)

if ($hpParagraph.Range.Text.TrimEnd() -ne "+2 Hit Points") {
    throw "Unexpected paragraph after 'Paladin' heading: '$($hpParagraph.Range.Text)'"
}

# Delete the whole paragraph, including its trailing paragraph mark, so the
# bullet list closes back up (Paladin -> 1 Color Slot -> Create White -> ...).
$hpParagraph.Range.Delete() | Out-Null

Write-Output "Removed the Paladin '+2 Hit Points' bullet."
